# Insert a new data row at row 200 (pushing existing rows 200:266 down to
# 201:267) and populate it with a new "Choclero" price record for
# Región Metropolitana dated 2022-04-12 (serial 44663).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(200).EntireRow.Insert()

$ws.Cells.Item(200, 1).Value  = 4
$ws.Cells.Item(200, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(200, 3).Value  = "Los Lagos"
$ws.Cells.Item(200, 4).Value  = 44663
$ws.Cells.Item(200, 5).Value  = 10
$ws.Cells.Item(200, 6).Value  = 100112024
$ws.Cells.Item(200, 7).Value  = "Choclo"
$ws.Cells.Item(200, 8).Value  = "Choclero"
$ws.Cells.Item(200, 9).Value  = "Primera"
$ws.Cells.Item(200, 10).Value = 6000
$ws.Cells.Item(200, 11).Value = 300
$ws.Cells.Item(200, 12).Value = 400
$ws.Cells.Item(200, 13).Value = 350
$ws.Cells.Item(200, 14).Value = "`$/unidad"
$ws.Cells.Item(200, 15).Value = "Región Metropolitana"
$ws.Cells.Item(200, 16).Value = 350
$ws.Cells.Item(200, 17).Value = 1
$ws.Cells.Item(200, 18).Value = "Hortaliza"
